$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 201
$ws.Range("I9").Value = 212.08333
$ws.Range("J9").Value = 156.66667
$ws.Range("K9").Value = 212.08333
$ws.Range("L9").Value = 156.66667
$ws.Range("M9").Value = -43.08332999999999
$ws.Range("N9").Value = -494.66667
$ws.Range("H43").Value = 9555.083000000001
$ws.Range("J43").Value = 5631
$ws.Range("L43").Value = 5631
$ws.Range("N43").Value = -5769
$ws.Range("H57").Value = 43260
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 43260
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 129780
$ws.Range("M57").Value = ""
$ws.Range("N57").Value = -130778
$ws.Range("H70").Value = 8838.519
$ws.Range("J70").Value = 7005.222
$ws.Range("L70").Value = 21015.666
$ws.Range("N70").Value = -21555.666
$ws.Range("H73").Value = 8838.519
$ws.Range("J73").Value = 7005.222
$ws.Range("L73").Value = 21015.666
$ws.Range("N73").Value = -22887.666
$ws.Range("H76").Value = 4501.7144
$ws.Range("I76").Value = 4013.0715
$ws.Range("J76").Value = 5479
$ws.Range("K76").Value = 4013.0715
$ws.Range("L76").Value = 5479
$ws.Range("M76").Value = -3698.0715
$ws.Range("N76").Value = -6109
$ws.Range("H79").Value = 4501.7144
$ws.Range("I79").Value = 4013.0715
$ws.Range("J79").Value = 5479
$ws.Range("K79").Value = 4013.0715
$ws.Range("L79").Value = 5479
$ws.Range("M79").Value = -2921.0715
$ws.Range("N79").Value = -7663
$ws.Range("H86").Value = 2848.4211
$ws.Range("I86").Value = 2834.2222
$ws.Range("K86").Value = 2834.2222
$ws.Range("M86").Value = -1711.2222
$ws.Range("H89").Value = 2848.4211
$ws.Range("I89").Value = 2834.2222
$ws.Range("K89").Value = 14171.111
$ws.Range("M89").Value = -8555.111000000001
$ws.Range("H112").Value = 7513.41
$ws.Range("I112").Value = 1179
$ws.Range("J112").Value = 8041.278
$ws.Range("K112").Value = 3537
$ws.Range("L112").Value = 24123.834
$ws.Range("M112").Value = -2429
$ws.Range("N112").Value = -26339.834
$ws.Range("H128").Value = 34166.668
$ws.Range("J128").Value = 34166.668
$ws.Range("L128").Value = 34166.668
$ws.Range("N128").Value = -44126.668
$ws.Range("H129").Value = 855.8
$ws.Range("I129").Value = 692.75
$ws.Range("J129").Value = 1508
$ws.Range("K129").Value = 2078.25
$ws.Range("L129").Value = 4524
$ws.Range("M129").Value = 2921.75
$ws.Range("N129").Value = -14524
$ws.Range("H135").Value = 15603.895
$ws.Range("J135").Value = 49803.273
$ws.Range("L135").Value = 448229.457
$ws.Range("N135").Value = -453299.457
$ws.Range("H137").Value = 8227.741
$ws.Range("I137").Value = 16420.238
$ws.Range("K137").Value = 49260.71400000001
$ws.Range("M137").Value = -46710.71400000001
$ws.Range("H138").Value = 1629.3846
$ws.Range("I138").Value = 1629.3846
$ws.Range("K138").Value = 4888.1538
$ws.Range("M138").Value = 251.8462
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2330.4167
$ws.Range("I2").Value = 2415.7778
$ws.Range("K2").Value = 2415.7778
$ws.Range("M2").Value = -2302.7778
$ws.Range("H63").Value = 4043.3333
$ws.Range("H66").Value = 4043.3333
$ws.Range("H74").Value = 667991.1
$ws.Range("I74").Value = 667991.1
$ws.Range("K74").Value = 667991.1
$ws.Range("M74").Value = -667117.1
$ws.Range("H77").Value = 667991.1
$ws.Range("I77").Value = 667991.1
$ws.Range("K77").Value = 3339955.5
$ws.Range("M77").Value = -3335587.5
$ws.Range("H92").Value = 20045000
$ws.Range("J92").Value = 40000000
$ws.Range("L92").Value = 40000000
$ws.Range("N92").Value = -40004992
$ws.Range("H102").Value = 5361.1333
$ws.Range("I102").Value = 4933.227
$ws.Range("K102").Value = 4933.227
$ws.Range("M102").Value = -3311.227
$ws.Range("H116").Value = 2330.4167
$ws.Range("I116").Value = 2415.7778
$ws.Range("K116").Value = 2415.7778
$ws.Range("M116").Value = -121.7777999999998
$ws.Range("H122").Value = 3611.8
$ws.Range("I122").Value = 3476
$ws.Range("J122").Value = 3747.6
$ws.Range("K122").Value = 10428
$ws.Range("L122").Value = 11242.8
$ws.Range("M122").Value = -7978
$ws.Range("N122").Value = -16142.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2330.4167
$ws.Range("I3").Value = 2415.7778
$ws.Range("K3").Value = 2415.7778
$ws.Range("M3").Value = -2301.7778
$ws.Range("H8").Value = 9667.166999999999
$ws.Range("I8").Value = 9400.6
$ws.Range("K8").Value = 9400.6
$ws.Range("M8").Value = -9260.6
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H88").Value = 67500
$ws.Range("J88").Value = 67500
$ws.Range("L88").Value = 67500
$ws.Range("N88").Value = -68312
$ws.Range("H91").Value = 67500
$ws.Range("J91").Value = 67500
$ws.Range("L91").Value = 67500
$ws.Range("N91").Value = -70308
$ws.Range("H99").Value = 2223.6365
$ws.Range("I99").Value = 983.8261
$ws.Range("J99").Value = 5075.2
$ws.Range("K99").Value = 983.8261
$ws.Range("L99").Value = 5075.2
$ws.Range("M99").Value = 514.1739
$ws.Range("N99").Value = -8071.2
$ws.Range("H105").Value = 4398.885
$ws.Range("I105").Value = 2793.6
$ws.Range("J105").Value = 9749.833000000001
$ws.Range("K105").Value = 2793.6
$ws.Range("L105").Value = 9749.833000000001
$ws.Range("M105").Value = -1046.6
$ws.Range("N105").Value = -13243.833
$ws.Range("H107").Value = 33214.812
$ws.Range("I107").Value = 57539
$ws.Range("J107").Value = 1940.8572
$ws.Range("K107").Value = 57539
$ws.Range("L107").Value = 1940.8572
$ws.Range("M107").Value = -55619
$ws.Range("N107").Value = -5780.8572
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2859298.8
$ws.Range("I31").Value = 2859298.8
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2859298.8
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2859003.8
$ws.Range("N31").Value = ""
$ws.Range("H34").Value = 2859298.8
$ws.Range("I34").Value = 2859298.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2859298.8
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2859096.8
$ws.Range("N34").Value = ""
$ws.Range("H62").Value = 6040.1763
$ws.Range("I62").Value = 2711.5715
$ws.Range("J62").Value = 8370.200000000001
$ws.Range("K62").Value = 2711.5715
$ws.Range("L62").Value = 8370.200000000001
$ws.Range("M62").Value = -2087.5715
$ws.Range("N62").Value = -9618.200000000001
$ws.Range("H65").Value = 6040.1763
$ws.Range("I65").Value = 2711.5715
$ws.Range("J65").Value = 8370.200000000001
$ws.Range("K65").Value = 13557.8575
$ws.Range("L65").Value = 41851
$ws.Range("M65").Value = -10437.8575
$ws.Range("N65").Value = -48091
$ws.Range("H94").Value = 1575.1666
$ws.Range("J94").Value = 1092.2
$ws.Range("L94").Value = 1092.2
$ws.Range("N94").Value = -1994.2
$ws.Range("H105").Value = 2067.8667
$ws.Range("I105").Value = 1136.5454
$ws.Range("J105").Value = 4629
$ws.Range("K105").Value = 1136.5454
$ws.Range("L105").Value = 4629
$ws.Range("M105").Value = 610.4546
$ws.Range("N105").Value = -8123
$ws.Range("H107").Value = 970.4167
$ws.Range("I107").Value = 779.2
$ws.Range("J107").Value = 1107
$ws.Range("K107").Value = 779.2
$ws.Range("L107").Value = 1107
$ws.Range("M107").Value = 1140.8
$ws.Range("N107").Value = -4947
$ws.Range("H132").Value = 16948.13
$ws.Range("I132").Value = 18671.53
$ws.Range("K132").Value = 56014.59
$ws.Range("M132").Value = -53484.59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3846245.5
$ws.Range("J7").Value = 19230784
$ws.Range("L7").Value = 57692352
$ws.Range("N7").Value = -57692576
$ws.Range("H12").Value = 217.22223
$ws.Range("J12").Value = 192.08333
$ws.Range("L12").Value = 576.24999
$ws.Range("N12").Value = -922.24999
$ws.Range("H17").Value = 222.4
$ws.Range("I17").Value = 226.66667
$ws.Range("K17").Value = 680.00001
$ws.Range("M17").Value = -511.00001
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H50").Value = 76.833336
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 76.833336
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 230.500008
$ws.Range("M50").Value = ""
$ws.Range("N50").Value = -1192.500008
$ws.Range("H53").Value = 76.833336
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 76.833336
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 230.500008
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = -1192.500008
$ws.Range("H68").Value = 4636
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622
$ws.Range("H71").Value = 4636
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112
$ws.Range("H116").Value = 4530.222
$ws.Range("I116").Value = 3038.8572
$ws.Range("K116").Value = 9116.571599999999
$ws.Range("M116").Value = -5674.571599999999
$ws.Range("H131").Value = 119862.39
$ws.Range("J131").Value = 2756.8125
$ws.Range("L131").Value = 8270.4375
$ws.Range("N131").Value = -18350.4375
$ws.Range("H137").Value = 3068.7
$ws.Range("J137").Value = 2748.5
$ws.Range("L137").Value = 8245.5
$ws.Range("N137").Value = -18445.5
$ws.Range("H141").Value = 5400.364
$ws.Range("I141").Value = 6201.143
$ws.Range("K141").Value = 18603.429
$ws.Range("M141").Value = -13423.429
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 504000
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -999884
$ws.Range("N3").Value = -8232
$ws.Range("H80").Value = 9709.799999999999
$ws.Range("I80").Value = 4933
$ws.Range("K80").Value = 4933
$ws.Range("M80").Value = -3935
$ws.Range("H83").Value = 9709.799999999999
$ws.Range("I83").Value = 4933
$ws.Range("K83").Value = 24665
$ws.Range("M83").Value = -19673
$ws.Range("H102").Value = 23237.13
$ws.Range("I102").Value = 27447.895
$ws.Range("J102").Value = 3236
$ws.Range("K102").Value = 27447.895
$ws.Range("L102").Value = 3236
$ws.Range("M102").Value = -25825.895
$ws.Range("N102").Value = -6480
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""
$ws.Range("H122").Value = 4399.4
$ws.Range("I122").Value = 4285
$ws.Range("J122").Value = 4666.3335
$ws.Range("K122").Value = 12855
$ws.Range("L122").Value = 13999.0005
$ws.Range("M122").Value = -10405
$ws.Range("N122").Value = -18899.0005
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("H132").Value = 1186.5
$ws.Range("I132").Value = 839.8182
$ws.Range("K132").Value = 2519.4546
$ws.Range("M132").Value = 10.54539999999997
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1655.0476
$ws.Range("I55").Value = 659.8
$ws.Range("J55").Value = 1966.0625
$ws.Range("K55").Value = 659.8
$ws.Range("L55").Value = 1966.0625
$ws.Range("M55").Value = -486.8
$ws.Range("N55").Value = -2312.0625
$ws.Range("H100").Value = 2836.125
$ws.Range("I100").Value = 2749.75
$ws.Range("K100").Value = 2749.75
$ws.Range("M100").Value = -2208.75
$ws.Range("I132").Value = 3732.5557
$ws.Range("J132").Value = 5684.857
$ws.Range("K132").Value = 11197.6671
$ws.Range("L132").Value = 17054.571
$ws.Range("M132").Value = -8667.667099999999
$ws.Range("N132").Value = -22114.571
$ws.Range("H134").Value = 124981.25
$ws.Range("J134").Value = 124984.5
$ws.Range("L134").Value = 124984.5
$ws.Range("N134").Value = -135124.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 114734.6
$ws.Range("J46").Value = 114734.6
$ws.Range("L46").Value = 114734.6
$ws.Range("N46").Value = -115196.6
$ws.Range("H81").Value = 7395.2
$ws.Range("I81").Value = 7395.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 14790.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -13729.4
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 7395.2
$ws.Range("I84").Value = 7395.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 73952
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -68648
$ws.Range("N84").Value = ""
$ws.Range("H109").Value = 87499.5
$ws.Range("J109").Value = 87499.5
$ws.Range("L109").Value = 87499.5
$ws.Range("N109").Value = -90273.5
$ws.Range("H132").Value = 23577.656
$ws.Range("I132").Value = 26972.5
$ws.Range("K132").Value = 80917.5
$ws.Range("M132").Value = -78387.5
$ws.Range("H134").Value = 114734.6
$ws.Range("J134").Value = 114734.6
$ws.Range("L134").Value = 344203.8
$ws.Range("N134").Value = -349273.8
$ws.Range("H140").Value = 92500.5
$ws.Range("J140").Value = 92500.5
$ws.Range("L140").Value = 92500.5
$ws.Range("N140").Value = -102860.5
